# Change the East-Asian font used by the "Normal" and "Heading" paragraph
# styles from "DejaVu Sans" to "Tahoma" (w:rFonts/@w:eastAsia).
$d = $word.ActiveDocument

$normal = $d.Styles.Item("Normal")
$normal.Font.NameFarEast = "Tahoma"

$heading = $d.Styles.Item("Heading")
$heading.Font.NameFarEast = "Tahoma"

# The "List", "Caption" and "Index" styles previously had no explicit
# complex-script font; give them an explicit w:rFonts/@w:cs="DejaVu Sans"
# (Font.NameBi maps to w:cs) matching the rest of the theme.
$list = $d.Styles.Item("List")
$list.Font.NameBi = "DejaVu Sans"

$caption = $d.Styles.Item("Caption")
$caption.Font.NameBi = "DejaVu Sans"

$index = $d.Styles.Item("Index")
$index.Font.NameBi = "DejaVu Sans"
